$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("1er Parcial")
$ws.Range("I4").Value = 8.7
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 35.48
$ws.Range("I5").Value = 7.8
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 33.33
$ws.Range("I11").Value = 6.8
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 4.35
$ws.Range("E12").Value = 24
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = 61.54
$ws.Range("H12").Value = 38.46
$ws.Range("I12").Value = 8
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 38.46

$ws = $wb.Worksheets.Item("2o Parcial")
$ws.Range("E2").Value = 29
$ws.Range("F2").Value = 11
$ws.Range("G2").Value = 72.5
$ws.Range("H2").Value = 27.5
$ws.Range("I2").Value = 7.8
$ws.Range("J2").Value = 11
$ws.Range("K2").Value = 27.5
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = 13.89
$ws.Range("H3").Value = 86.11
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 31
$ws.Range("K3").Value = 86.11
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 18
$ws.Range("G4").Value = 41.94
$ws.Range("H4").Value = 58.06
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 18
$ws.Range("K4").Value = 58.06
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 42.86
$ws.Range("H5").Value = 57.14
$ws.Range("I5").Value = 8.3
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 57.14
$ws.Range("E6").Value = 26
$ws.Range("F6").Value = 15
$ws.Range("G6").Value = 63.41
$ws.Range("H6").Value = 36.59
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 15
$ws.Range("K6").Value = 36.59
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = 18
$ws.Range("G10").Value = 50
$ws.Range("H10").Value = 50
$ws.Range("I10").Value = 8.3
$ws.Range("J10").Value = 18
$ws.Range("K10").Value = 50
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 69.57
$ws.Range("H11").Value = 30.43
$ws.Range("I11").Value = 6.6
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 4.35
$ws.Range("E12").Value = 22
$ws.Range("F12").Value = 17
$ws.Range("G12").Value = 56.41
$ws.Range("H12").Value = 43.59
$ws.Range("I12").Value = 8.1
$ws.Range("J12").Value = 17
$ws.Range("K12").Value = 43.59

$ws = $wb.Worksheets.Item("3er Parcial")
$ws.Range("I3").Value = 8
$ws.Range("E4").Value = 20
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = 64.52
$ws.Range("H4").Value = 35.48
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 11
$ws.Range("K4").Value = 35.48
$ws.Range("E5").Value = 14
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = 66.67
$ws.Range("H5").Value = 33.33
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 33.33
$ws.Range("I6").Value = 8
$ws.Range("I10").Value = 8.4
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 69.57
$ws.Range("H11").Value = 30.43
$ws.Range("I11").Value = 6.7
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 4.35
$ws.Range("E12").Value = 24
$ws.Range("F12").Value = 15
$ws.Range("G12").Value = 61.54
$ws.Range("H12").Value = 38.46
$ws.Range("J12").Value = 15
$ws.Range("K12").Value = 38.46
